# ADICIONADO INSERÇÃO DE CHV FACA
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Replace the numeric "chave" in A2 with the text value (keeps existing style s="1")
$ws.Range("A2").Value = "1145623903 A"

# Insert the new "chave" text value in A3, matching the style already used by J3 (s="2": centered horizontally)
$ws.Range("A3").Value = "1145623903 B"
$ws.Range("A3").HorizontalAlignment = -4108  # xlCenter

# Move the active selection to A4, as it is after the insertion
$ws.Range("A4").Select()

$ws.Activate()
